$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price cells so Excel keeps them as text (matches source inlineStr cells)
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D18", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '30.055.52'
$ws.Range("E2").Value = '  +0.63%  '
$ws.Range("D3").Value = '1.889.89'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("D4").Value = '0.9978'
$ws.Range("E4").Value = '  -0.35%  '
$ws.Range("D5").Value = '0.7442'
$ws.Range("E5").Value = '  -1.78%  '
$ws.Range("D6").Value = '242.84'
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("D7").Value = '0.9975'
$ws.Range("E7").Value = '  -0.42%  '
$ws.Range("D8").Value = '0.3165'
$ws.Range("E8").Value = '  +1.42%  '
$ws.Range("D9").Value = '0.07246'
$ws.Range("E9").Value = '  +1.39%  '
$ws.Range("D10").Value = '24.99'
$ws.Range("E10").Value = '  -1.61%  '
$ws.Range("D11").Value = '0.08361'
$ws.Range("E11").Value = '  -2.02%  '
$ws.Range("D12").Value = '0.7606'
$ws.Range("E12").Value = '  -0.05%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.912.77'
$ws.Range("E13").Value = '  +0.97%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '5.427'
$ws.Range("E14").Value = '  +1.28%  '
$ws.Range("D15").Value = '92.74'
$ws.Range("E15").Value = '  -0.76%  '
$ws.Range("D16").Value = '6.154'
$ws.Range("E16").Value = '  +0.35%  '
$ws.Range("D17").Value = '30.037.25'
$ws.Range("E17").Value = '  +0.46%  '
$ws.Range("D18").Value = '249.64'
$ws.Range("E18").Value = '  +2.37%  '
$ws.Range("E19").Value = '  -0.55%  '
$ws.Range("D20").Value = '0.000007869'
$ws.Range("E20").Value = '  +0.77%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '0.9964'
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.132.79'
$ws.Range("E22").Value = '  -1.83%  '
$ws.Range("D23").Value = '8.009'
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").Value = '0.9981'
$ws.Range("E24").Value = '  -0.27%  '
$ws.Range("D25").Value = '0.1579'
$ws.Range("E25").Value = '  -1.90%  '
$ws.Range("D26").Value = '9.306'
$ws.Range("E26").Value = '  -0.73%  '
$ws.Range("D27").Value = '164.67'
$ws.Range("E27").Value = '  +1.06%  '
$ws.Range("D28").Value = '18.75'
$ws.Range("D29").Value = '2.054'
$ws.Range("E29").Value = '  +1.41%  '
$ws.Range("D30").Value = '1.483'
$ws.Range("E30").Value = '  -2.48%  '
$ws.Range("D31").Value = '4.605'
$ws.Range("E31").Value = '  +2.83%  '
$ws.Range("D32").Value = '1.537'
$ws.Range("E32").Value = '  +0.30%  '
$ws.Range("D33").Value = '4.223'
$ws.Range("E33").Value = '  +2.86%  '
$ws.Range("D34").Value = '0.05379'
$ws.Range("E34").Value = '  -0.71%  '
$ws.Range("D35").Value = '1.254'
$ws.Range("E35").Value = '  +1.16%  '
$ws.Range("D36").Value = '0.7616'
$ws.Range("E36").Value = '  +2.37%  '
$ws.Range("D37").Value = '0.9958'
$ws.Range("E37").Value = '  -0.47%  '
$ws.Range("D38").Value = '2.720'
$ws.Range("E38").Value = '  +0.40%  '
$ws.Range("D39").Value = '0.01972'
$ws.Range("E39").Value = '  +1.67%  '
$ws.Range("E40").Value = '  -0.33%  '
$ws.Range("D41").Value = '0.4572'
$ws.Range("E41").Value = '  +2.64%  '
$ws.Range("D42").Value = '1.100.68'
$ws.Range("E42").Value = '  -0.29%  '
$ws.Range("D43").Value = '73.00'
$ws.Range("E43").Value = '  +0.60%  '
$ws.Range("D44").Value = '6.055'
$ws.Range("E44").Value = '  -0.35%  '
$ws.Range("D45").Value = '0.8735'
$ws.Range("E45").Value = '  +1.59%  '
$ws.Range("D46").Value = '104.59'
$ws.Range("E46").Value = '  +1.90%  '
$ws.Range("D47").Value = '0.9985'
$ws.Range("E47").Value = '  -0.29%  '
$ws.Range("D48").Value = '1.874'
$ws.Range("E49").Value = '  -0.27%  '
$ws.Range("D50").Value = '9.604'
$ws.Range("E50").Value = '  -1.17%  '
$ws.Range("D51").Value = '2.033.56'
$ws.Range("E51").Value = '  -0.78%  '
